# Update the cryptos report: refresh Price (D) and Volume(1h) (E) figures,
# and apply the Fetch.AI/ThetaToken and VeChain/WEMIXToken row swaps
# (rows 44-47), matching the latest GitHub Actions data refresh.
#
# Note: the Price column stores numeric-looking text (e.g. "181.06",
# "1.00", "0.0000254") as literal strings in the source file. Assigning
# such strings directly to Range.Value lets Excel auto-coerce them to
# numbers, which silently drops trailing zeros / introduces float
# rounding. To keep the exact original text semantics, those cells are
# temporarily switched to Text format ("@") before the write, then
# restored to the default "Normal" style afterwards so no stray
# cell-level formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.252.06'
$ws.Range("E2").Value = '  -6.59%  '
$ws.Range("D3").Value = '3.280.21'
$ws.Range("E3").Value = '  -7.53%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '181.06'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -11.57%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '518.46'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.85%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.595'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.23%  '
$ws.Range("D8").Value = '3.277.87'
$ws.Range("E8").Value = '  -7.44%  '
$ws.Range("E9").Value = '  -0.08%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.617'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -7.31%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.38'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -8.08%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.131'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -9.30%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000254'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -7.64%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.09'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -8.84%  '
$ws.Range("D15").Value = '3.808.93'
$ws.Range("E15").Value = '  -8.09%  '
$ws.Range("E16").Value = '  -5.78%  '
$ws.Range("D17").Value = '3.293.96'
$ws.Range("E17").Value = '  -7.97%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.65'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -6.44%  '
$ws.Range("D19").Value = '63.338.07'
$ws.Range("E19").Value = '  -6.83%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.92'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -8.98%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.946'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -9.76%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '370.04'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.85%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.22'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -8.23%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.68'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -9.99%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '79.85'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.23%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.82'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("E27").Value = '  -2.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.61'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -7.65%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '11.32'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -7.43%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.28'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.87%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '646.18'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -7.99%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.42'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -8.36%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.65'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -10.28%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.12'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.56%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.105'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '58.91'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -6.95%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.389'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.20%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.98'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -11.58%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.08%  '
$ws.Range("D41").Value = '2.986.35'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.125'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.37%  '
$ws.Range("D43").Value = '0.0₃0648'
$ws.Range("E43").Value = '  -10.05%  '
$ws.Range("B44").Value = 'Fetch.AI'
$ws.Range("C44").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.42'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.08%  '
$ws.Range("B45").Value = 'ThetaToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.67'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -15.73%  '
$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.60'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.96%  '
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0387'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.59%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.80'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.17%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.124'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.47'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -20.22%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.90'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.27%  '
